$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated capital structure figures to rows 2 and 3.
foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = 0.04269999999999999
    $ws.Range("E$row").ClearContents()

    $ws.Range("G$row").Value = 0.1173210161662818
    $ws.Range("H$row").Value = 0.1173210161662818
    $ws.Range("I$row").Value = 0.1053117782909931
    $ws.Range("J$row").Value = 0.06420874840577712
    $ws.Range("K$row").Value = 7.86
    $ws.Range("L$row").Value = 0.03630484988452656

    $ws.Range("U$row").Value = 52.9
    $ws.Range("V$row").Value = 0.2454756380510441
    $ws.Range("W$row").Value = 0.04316309719934103
    $ws.Range("X$row").Value = 0.04079568500000472
    $ws.Range("Y$row").Value = 0.002367412199336311
    $ws.Range("Z$row").Value = 1.523039043264158
    $ws.Range("AA$row").Value = 0.09779243074112379
    $ws.Range("AB$row").Value = 0.04009472595480686
    $ws.Range("AC$row").Value = 0.05769770478631693
    $ws.Range("AD$row").Value = 5.75
    $ws.Range("AE$row").Value = 0
    $ws.Range("AF$row").Value = 5.75
    $ws.Range("AG$row").Value = -47.15
    $ws.Range("AH$row").Value = 0.02598870056497175
    $ws.Range("AI$row").Value = 0.02444208289054198
    $ws.Range("AJ$row").Value = -0.2800712800712801
    $ws.Range("AK$row").Value = -0.2585686865917192
    $ws.Range("AL$row").Value = 0.744
    $ws.Range("AM$row").Value = 0.744
    $ws.Range("AN$row").Value = 0.2263779527559055
    $ws.Range("AO$row").Value = 30.64516129032258
    $ws.Range("AP$row").Value = -1.856299212598425
    $ws.Range("AQ$row").Value = 30.64516129032258
}
